$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2943.2
$ws.Range("I15").Value = 2943.2
$ws.Range("K15").Value = 8829.599999999999
$ws.Range("M15").Value = -8660.599999999999
$ws.Range("H17").Value = 25941.012
$ws.Range("J17").Value = 25941.012
$ws.Range("L17").Value = 77823.03599999999
$ws.Range("N17").Value = -78159.03599999999
$ws.Range("H100").Value = 2762.8635
$ws.Range("I100").Value = 1939.2858
$ws.Range("J100").Value = 4204.125
$ws.Range("K100").Value = 1939.2858
$ws.Range("L100").Value = 4204.125
$ws.Range("M100").Value = -1398.2858
$ws.Range("N100").Value = -5286.125
$ws.Range("H112").Value = 10459412
$ws.Range("I112").Value = 633.3333
$ws.Range("J112").Value = 11953523
$ws.Range("K112").Value = 1899.9999
$ws.Range("L112").Value = 35860569
$ws.Range("M112").Value = -791.9999
$ws.Range("N112").Value = -35862785
$ws.Range("H116").Value = 4873.16
$ws.Range("I116").Value = 4921.643
$ws.Range("J116").Value = 4811.4546
$ws.Range("K116").Value = 4921.643
$ws.Range("L116").Value = 4811.4546
$ws.Range("M116").Value = -1479.643
$ws.Range("N116").Value = -11695.4546
$ws.Range("H132").Value = 3706246.8
$ws.Range("I132").Value = 4880556
$ws.Range("J132").Value = 2656.5386
$ws.Range("K132").Value = 14641668
$ws.Range("L132").Value = 7969.6158
$ws.Range("M132").Value = -14639138
$ws.Range("N132").Value = -13029.6158
$ws.Range("H137").Value = 2844.8596
$ws.Range("I137").Value = 2981.6086
$ws.Range("J137").Value = 2273
$ws.Range("K137").Value = 8944.825800000001
$ws.Range("L137").Value = 6819
$ws.Range("M137").Value = -6394.825800000001
$ws.Range("N137").Value = -11919

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1599.9143
$ws.Range("I61").Value = 916.1724
$ws.Range("J61").Value = 4904.6665
$ws.Range("K61").Value = 916.1724
$ws.Range("L61").Value = 4904.6665
$ws.Range("M61").Value = -704.1724
$ws.Range("N61").Value = -5328.6665
$ws.Range("H74").Value = 706.4
$ws.Range("I74").Value = 581.2
$ws.Range("K74").Value = 581.2
$ws.Range("M74").Value = 292.8
$ws.Range("H77").Value = 706.4
$ws.Range("I77").Value = 581.2
$ws.Range("K77").Value = 2906
$ws.Range("M77").Value = 1462
$ws.Range("H97").Value = 999.8
$ws.Range("I97").Value = 999.8
$ws.Range("K97").Value = 999.8
$ws.Range("M97").Value = -503.8
$ws.Range("H102").Value = 3518.889
$ws.Range("I102").Value = 2708.75
$ws.Range("K102").Value = 2708.75
$ws.Range("M102").Value = -1086.75
$ws.Range("H136").Value = 1599.9143
$ws.Range("I136").Value = 916.1724
$ws.Range("J136").Value = 4904.6665
$ws.Range("K136").Value = 2748.5172
$ws.Range("L136").Value = 14713.9995
$ws.Range("M136").Value = -198.5172000000002
$ws.Range("N136").Value = -19813.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1240
$ws.Range("I94").Value = 1083.3334
$ws.Range("J94").Value = 1396.6666
$ws.Range("K94").Value = 1083.3334
$ws.Range("L94").Value = 1396.6666
$ws.Range("M94").Value = -632.3334
$ws.Range("N94").Value = -2298.6666
$ws.Range("H99").Value = 3036.5454
$ws.Range("I99").Value = 1505
$ws.Range("J99").Value = 3911.7144
$ws.Range("K99").Value = 1505
$ws.Range("L99").Value = 3911.7144
$ws.Range("M99").Value = -7
$ws.Range("N99").Value = -6907.7144
$ws.Range("H134").Value = 2301.6191
$ws.Range("I134").Value = 1449.0588
$ws.Range("K134").Value = 4347.1764
$ws.Range("M134").Value = -1812.1764

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1880.8
$ws.Range("I16").Value = 1416.5
$ws.Range("J16").Value = 2190.3333
$ws.Range("K16").Value = 1416.5
$ws.Range("L16").Value = 2190.3333
$ws.Range("M16").Value = -1129.5
$ws.Range("N16").Value = -2764.3333
$ws.Range("H31").Value = 3601.2778
$ws.Range("I31").Value = 2220.2727
$ws.Range("K31").Value = 2220.2727
$ws.Range("M31").Value = -1925.2727
$ws.Range("H34").Value = 3601.2778
$ws.Range("I34").Value = 2220.2727
$ws.Range("K34").Value = 2220.2727
$ws.Range("M34").Value = -2018.2727
$ws.Range("H58").Value = 7354969
$ws.Range("I58").Value = 1256.6604
$ws.Range("J58").Value = 33338086
$ws.Range("K58").Value = 1256.6604
$ws.Range("L58").Value = 33338086
$ws.Range("M58").Value = -1053.6604
$ws.Range("N58").Value = -33338492
$ws.Range("H113").Value = 1880.8
$ws.Range("I113").Value = 1416.5
$ws.Range("J113").Value = 2190.3333
$ws.Range("K113").Value = 1416.5
$ws.Range("L113").Value = 2190.3333
$ws.Range("M113").Value = 753.5
$ws.Range("N113").Value = -6530.3333
$ws.Range("H134").Value = 3255.7058
$ws.Range("I134").Value = 1695.5834
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 5086.7502
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -2551.7502
$ws.Range("N134").Value = -26070
$ws.Range("H136").Value = 7354969
$ws.Range("I136").Value = 1256.6604
$ws.Range("J136").Value = 33338086
$ws.Range("K136").Value = 3769.9812
$ws.Range("L136").Value = 100014258
$ws.Range("M136").Value = -1219.9812
$ws.Range("N136").Value = -100019358

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 2023.5
$ws.Range("I92").Value = 91
$ws.Range("J92").Value = 2506.625
$ws.Range("K92").Value = 273
$ws.Range("L92").Value = 7519.875
$ws.Range("M92").Value = 975
$ws.Range("N92").Value = -10015.875
$ws.Range("H129").Value = 1725.5238
$ws.Range("I129").Value = 1337.2727
$ws.Range("J129").Value = 2152.6
$ws.Range("K129").Value = 4011.8181
$ws.Range("L129").Value = 6457.799999999999
$ws.Range("M129").Value = 988.1819
$ws.Range("N129").Value = -16457.8
$ws.Range("H130").Value = 2171.6667
$ws.Range("I130").Value = 1030
$ws.Range("K130").Value = 3090
$ws.Range("M130").Value = 1930
$ws.Range("H131").Value = 1472.3243
$ws.Range("I131").Value = 1697.5
$ws.Range("J131").Value = 1410.2069
$ws.Range("K131").Value = 5092.5
$ws.Range("L131").Value = 4230.620699999999
$ws.Range("M131").Value = -52.5
$ws.Range("N131").Value = -14310.6207
$ws.Range("H134").Value = 2070.4
$ws.Range("I134").Value = 1096.1818
$ws.Range("J134").Value = 4749.5
$ws.Range("K134").Value = 3288.5454
$ws.Range("L134").Value = 14248.5
$ws.Range("M134").Value = 1781.4546
$ws.Range("N134").Value = -24388.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 38666.668
$ws.Range("J75").Value = 38666.668
$ws.Range("L75").Value = 38666.668
$ws.Range("N75").Value = -40414.668
$ws.Range("H78").Value = 38666.668
$ws.Range("J78").Value = 38666.668
$ws.Range("L78").Value = 116000.004
$ws.Range("N78").Value = -124736.004
$ws.Range("H97").Value = 1538.4166
$ws.Range("I97").Value = 492.85715
$ws.Range("J97").Value = 3002.2
$ws.Range("K97").Value = 492.85715
$ws.Range("L97").Value = 3002.2
$ws.Range("M97").Value = 3.14285000000001
$ws.Range("N97").Value = -3994.2
$ws.Range("H107").Value = 817.86206
$ws.Range("I107").Value = 396
$ws.Range("J107").Value = 1115.6471
$ws.Range("K107").Value = 396
$ws.Range("L107").Value = 1115.6471
$ws.Range("M107").Value = 1524
$ws.Range("N107").Value = -4955.6471
$ws.Range("H122").Value = 3352.8125
$ws.Range("I122").Value = 2482.7778
$ws.Range("J122").Value = 4471.4287
$ws.Range("K122").Value = 7448.3334
$ws.Range("L122").Value = 13414.2861
$ws.Range("M122").Value = -4998.3334
$ws.Range("N122").Value = -18314.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 166671380
$ws.Range("I61").Value = 333334080
$ws.Range("J61").Value = 8668.333000000001
$ws.Range("K61").Value = 333334080
$ws.Range("L61").Value = 8668.333000000001
$ws.Range("M61").Value = -333333878
$ws.Range("N61").Value = -9072.333000000001
$ws.Range("H93").Value = 3200.1667
$ws.Range("I93").Value = 2857.1428
$ws.Range("J93").Value = 4400.75
$ws.Range("K93").Value = 2857.1428
$ws.Range("L93").Value = 4400.75
$ws.Range("M93").Value = -1609.1428
$ws.Range("N93").Value = -6896.75
$ws.Range("H100").Value = 2169.2666
$ws.Range("I100").Value = 1324.75
$ws.Range("K100").Value = 1324.75
$ws.Range("M100").Value = -783.75
$ws.Range("H113").Value = 166671380
$ws.Range("I113").Value = 333334080
$ws.Range("J113").Value = 8668.333000000001
$ws.Range("K113").Value = 333334080
$ws.Range("L113").Value = 8668.333000000001
$ws.Range("M113").Value = -333331910
$ws.Range("N113").Value = -13008.333
$ws.Range("H132").Value = 2235.205
$ws.Range("I132").Value = 1350.7931
$ws.Range("K132").Value = 4052.379300000001
$ws.Range("M132").Value = -1522.379300000001
$ws.Range("H136").Value = 2406.94
$ws.Range("I136").Value = 1308.675
$ws.Range("J136").Value = 6800
$ws.Range("K136").Value = 3926.025
$ws.Range("L136").Value = 20400
$ws.Range("M136").Value = -1376.025
$ws.Range("N136").Value = -25500
$ws.Range("H138").Value = 30000
$ws.Range("J138").Value = 30000
$ws.Range("L138").Value = 30000
$ws.Range("N138").Value = -40280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 968.4
$ws.Range("I96").Value = 967.5714
$ws.Range("J96").Value = 980
$ws.Range("K96").Value = 967.5714
$ws.Range("L96").Value = 980
$ws.Range("M96").Value = 405.4286
$ws.Range("N96").Value = -3726
$ws.Range("H132").Value = 9239.5
$ws.Range("I132").Value = 1820.4038
$ws.Range("J132").Value = 33351.562
$ws.Range("K132").Value = 5461.2114
$ws.Range("L132").Value = 100054.686
$ws.Range("M132").Value = -2931.2114
$ws.Range("N132").Value = -105114.686
$ws.Range("H136").Value = 1055.2222
$ws.Range("I136").Value = 512.25
$ws.Range("K136").Value = 1536.75
